$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing ddl object cell in place (C3): "btnHome" -> "ddlSelectName".
# Since the old shared string is referenced nowhere else, the engine reuses
# the same shared-string slot instead of appending a new one.
$ws.Range("C3").Value = "ddlSelectName"

# Change the Action / UserInput for row 3 (D3, E3); F3 ("Y") stays the same.
$ws.Range("D3").Value = "Select"
$ws.Range("E3").Value = "Harry Potter"

# D3 needs the same "quote prefix" cell style as C2 (style index 2 in the
# original file). Setting .Value resets the style, so re-apply the format
# afterwards via a formats-only paste from a cell that already has it.
$ws.Range("C2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add a new row (row 4) for the "btnLogin" click step.
$ws.Range("B4").Value = "LoginPage"
$ws.Range("C4").Value = "btnLogin"
$ws.Range("D4").Value = "Click"
$ws.Range("E4").Value = "NA"
$ws.Range("F4").Value = "Y"

# Update the sheet's active selection to the new last cell, F4.
$ws.Range("F4").Select()
